# Update the data table in sheet1 to reflect the latest model run, and
# remove the trailing rows (28, 29, 30 -> old row labels 30-32) that are
# no longer part of the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values (rows are keyed by worksheet row number) for columns C:F
# (C = Detected Infections, D = Cumulative Infections,
#  E = Current Asymptomatic Infections, F = Current Infectious Carriers).
$data = @{
    10 = @(0, 1, 0, 1)
    11 = @(0, 2, 1, 1)
    12 = @(0, 6, 5, 1)
    13 = @(1, 6, 5, 0)
    14 = @(1, 8, 7, 0)
    15 = @(1, 8, 7, 0)
    16 = @(1, 8, 6, 1)
    17 = @(1, 10, 4, 5)
    18 = @(1, 15, 9, 5)
    19 = @(1, 21, 13, 7)
    20 = @(1, 29, 21, 7)
    21 = @(2, 33, 25, 6)
    22 = @(6, 35, 26, 3)
    23 = @(6, 37, 25, 6)
    24 = @(8, 42, 22, 12)
    25 = @(8, 46, 20, 18)
    26 = @(8, 56, 25, 23)
    27 = @(9, 66, 32, 25)
    28 = @(12, 77, 42, 23)
    29 = @(20, 87, 46, 21)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 3).Value = $vals[0]
    $ws.Cells.Item($row, 4).Value = $vals[1]
    $ws.Cells.Item($row, 5).Value = $vals[2]
    $ws.Cells.Item($row, 6).Value = $vals[3]
}

# Remove the old trailing rows 30-32 (previously day 28-30), the dataset
# now ends at row 29 (day 27).
$ws.Range("A30:F32").Delete()
